# Updated symbol list on Sun Jan 22 08:11:07 UTC 2023 with GitHub Actions
# Refreshes Price (D), Volume(1h) (E), and Hora (G) columns for the crypto
# ranking table on Sheet1, cell by cell, matching the source scrape diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each target cell already holds a text value (prices/percentages/hour
# shown as strings). Force the cell to Text format before writing so
# Excel does not reinterpret the new value as a number (which would both
# change the stored type and round values like "37.41"), then clear the
# format back off so no stray style gets left behind.
function Set-TextValue($addr, $value) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue "E2" "-0.72%"
Set-TextValue "G2" "8"
Set-TextValue "D3" "37.41"
Set-TextValue "E3" "6.40%"
Set-TextValue "G3" "8"
Set-TextValue "D4" "5.007"
Set-TextValue "E4" "-2.76%"
Set-TextValue "G4" "8"
Set-TextValue "D5" "0.07837"
Set-TextValue "E5" "0.48%"
Set-TextValue "G5" "8"
Set-TextValue "D6" "2.193"
Set-TextValue "E6" "-8.16%"
Set-TextValue "G6" "8"
Set-TextValue "D7" "8.023"
Set-TextValue "E7" "-0.45%"
Set-TextValue "G7" "8"
Set-TextValue "D8" "4.021"
Set-TextValue "E8" "1.76%"
Set-TextValue "G8" "8"
Set-TextValue "D9" "0.9090"
Set-TextValue "E9" "-1.60%"
Set-TextValue "G9" "8"
Set-TextValue "D10" "0.09672"
Set-TextValue "E10" "-2.82%"
Set-TextValue "G10" "8"
Set-TextValue "D11" "0.1893"
Set-TextValue "E11" "4.37%"
Set-TextValue "G11" "8"
Set-TextValue "D12" "0.08526"
Set-TextValue "E12" "-1.62%"
Set-TextValue "G12" "8"
Set-TextValue "E13" "6.32%"
Set-TextValue "G13" "8"
Set-TextValue "D14" "0.09960"
Set-TextValue "E14" "0.47%"
Set-TextValue "G14" "8"
Set-TextValue "D15" "0.001485"
Set-TextValue "E15" "-0.92%"
Set-TextValue "G15" "8"
Set-TextValue "D16" "0.005669"
Set-TextValue "E16" "-0.46%"
Set-TextValue "G16" "8"
Set-TextValue "D17" "3.465"
Set-TextValue "E17" "-0.04%"
Set-TextValue "G17" "8"
Set-TextValue "D18" "2.070"
Set-TextValue "E18" "-4.32%"
Set-TextValue "G18" "8"
Set-TextValue "E19" "2.79%"
Set-TextValue "G19" "8"
Set-TextValue "D20" "0.1301"
Set-TextValue "E20" "-1.79%"
Set-TextValue "G20" "8"
Set-TextValue "E21" "9.44%"
Set-TextValue "G21" "8"
Set-TextValue "D22" "0.2204"
Set-TextValue "E22" "-7.57%"
Set-TextValue "G22" "8"
Set-TextValue "D23" "0.04631"
Set-TextValue "E23" "1.32%"
Set-TextValue "G23" "8"
Set-TextValue "D24" "0.001231"
Set-TextValue "E24" "1.07%"
Set-TextValue "G24" "8"
Set-TextValue "D25" "0.004802"
Set-TextValue "E25" "7.92%"
Set-TextValue "G25" "8"
Set-TextValue "E26" "-7.58%"
Set-TextValue "G26" "8"
Set-TextValue "D27" "0.0004754"
Set-TextValue "E27" "28.59%"
Set-TextValue "G27" "8"
Set-TextValue "G28" "8"
Set-TextValue "G29" "8"
Set-TextValue "G30" "8"
Set-TextValue "G31" "8"
Set-TextValue "G32" "8"
Set-TextValue "G33" "8"
Set-TextValue "G34" "8"
Set-TextValue "G35" "8"
Set-TextValue "G36" "8"
Set-TextValue "G37" "8"
Set-TextValue "G38" "8"
Set-TextValue "D39" "0.01756"
Set-TextValue "E39" "-2.12%"
Set-TextValue "G39" "8"
Set-TextValue "D40" "0.04719"
Set-TextValue "E40" "-1.06%"
Set-TextValue "G40" "8"
Set-TextValue "D41" "0.007893"
Set-TextValue "E41" "1.79%"
Set-TextValue "G41" "8"
Set-TextValue "D42" "0.1392"
Set-TextValue "E42" "-1.24%"
Set-TextValue "G42" "8"
Set-TextValue "D43" "0.007662"
Set-TextValue "E43" "7.02%"
Set-TextValue "G43" "8"
Set-TextValue "D44" "0.002172"
Set-TextValue "E44" "-6.33%"
Set-TextValue "G44" "8"
Set-TextValue "D45" "0.009870"
Set-TextValue "E45" "3.43%"
Set-TextValue "G45" "8"
Set-TextValue "D46" "0.00006082"
Set-TextValue "E46" "-0.54%"
Set-TextValue "G46" "8"
Set-TextValue "D47" "0.00000000751"
Set-TextValue "E47" "0.15%"
Set-TextValue "G47" "8"
Set-TextValue "G48" "8"
Set-TextValue "E49" "34.55%"
Set-TextValue "G49" "8"
Set-TextValue "D50" "0.00002102"
Set-TextValue "E50" "0.15%"
Set-TextValue "G50" "8"
Set-TextValue "D51" "0.0002001"
Set-TextValue "E51" "0.15%"
Set-TextValue "G51" "8"
